$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A37").Value = 4
$ws.Range("B37").Value = "2：27 - 5：23"
$ws.Range("C37").Value = "柯西中值"

$ws.Range("C37").Select()
